$wb = $excel.ActiveWorkbook

# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H69").Value = 4925.7144
$ws.Range("I69").Value = 3750
$ws.Range("J69").Value = 5396
$ws.Range("K69").Value = 11250
$ws.Range("L69").Value = 16188
$ws.Range("M69").Value = -10376
$ws.Range("N69").Value = -17936
$ws.Range("H72").Value = 4925.7144
$ws.Range("I72").Value = 3750
$ws.Range("J72").Value = 5396
$ws.Range("K72").Value = 33750
$ws.Range("L72").Value = 48564
$ws.Range("M72").Value = -29382
$ws.Range("N72").Value = -57300
$ws.Range("H88").Value = 974.13336
$ws.Range("J88").Value = 1276.6
$ws.Range("L88").Value = 1276.6
$ws.Range("N88").Value = -2088.6
$ws.Range("H91").Value = 974.13336
$ws.Range("J91").Value = 1276.6
$ws.Range("L91").Value = 1276.6
$ws.Range("N91").Value = -4084.6
$ws.Range("H98").Value = 3409.425
$ws.Range("I98").Value = 3010.7297
$ws.Range("J98").Value = 8326.666999999999
$ws.Range("K98").Value = 3010.7297
$ws.Range("L98").Value = 8326.666999999999
$ws.Range("M98").Value = -1512.7297
$ws.Range("N98").Value = -11322.667
$ws.Range("H99").Value = 728.5
$ws.Range("I99").Value = 504.66666
$ws.Range("K99").Value = 1513.99998
$ws.Range("M99").Value = -15.99998000000005
$ws.Range("H115").Value = 5195
$ws.Range("I115").Value = 6350.7144
$ws.Range("J115").Value = 1150
$ws.Range("K115").Value = 19052.1432
$ws.Range("L115").Value = 3450
$ws.Range("M115").Value = -17485.1432
$ws.Range("N115").Value = -6584
$ws.Range("H122").Value = 3409.425
$ws.Range("I122").Value = 3010.7297
$ws.Range("J122").Value = 8326.666999999999
$ws.Range("K122").Value = 9032.1891
$ws.Range("L122").Value = 24980.001
$ws.Range("M122").Value = -6582.1891
$ws.Range("N122").Value = -29880.001
$ws.Range("H125").Value = 1000.7241
$ws.Range("I125").Value = 670.9524
$ws.Range("J125").Value = 1866.375
$ws.Range("K125").Value = 6038.5716
$ws.Range("L125").Value = 16797.375
$ws.Range("M125").Value = -3578.5716
$ws.Range("N125").Value = -21717.375
$ws.Range("H127").Value = 1236.8422
$ws.Range("I127").Value = 525
$ws.Range("J127").Value = 1426.6666
$ws.Range("K127").Value = 1575
$ws.Range("L127").Value = 4279.9998
$ws.Range("M127").Value = 3385
$ws.Range("N127").Value = -14199.9998
$ws.Range("H131").Value = 4736.136
$ws.Range("I131").Value = 1899.5454
$ws.Range("J131").Value = 7572.727
$ws.Range("K131").Value = 5698.6362
$ws.Range("L131").Value = 22718.181
$ws.Range("M131").Value = -658.6361999999999
$ws.Range("N131").Value = -32798.181
$ws.Range("H138").Value = 2899.7556
$ws.Range("I138").Value = 2615.1538
$ws.Range("J138").Value = 4749.6665
$ws.Range("K138").Value = 7845.4614
$ws.Range("L138").Value = 14248.9995
$ws.Range("M138").Value = -2705.4614
$ws.Range("N138").Value = -24528.9995

# ----- Sheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H3").Value = 41669416
$ws.Range("I3").Value = 250001500
$ws.Range("K3").Value = 250001500
$ws.Range("M3").Value = -250001385
$ws.Range("H122").Value = 2295.3157
$ws.Range("I122").Value = 1702.4
$ws.Range("J122").Value = 2507.0715
$ws.Range("K122").Value = 5107.200000000001
$ws.Range("L122").Value = 7521.2145
$ws.Range("M122").Value = -2657.200000000001
$ws.Range("N122").Value = -12421.2145

# ----- Sheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3102.16
$ws.Range("I31").Value = 1510.8889
$ws.Range("K31").Value = 1510.8889
$ws.Range("M31").Value = -1215.8889
$ws.Range("H34").Value = 3102.16
$ws.Range("I34").Value = 1510.8889
$ws.Range("K34").Value = 1510.8889
$ws.Range("M34").Value = -1308.8889
$ws.Range("H99").Value = 144259.14
$ws.Range("I99").Value = 1400
$ws.Range("J99").Value = 334738
$ws.Range("K99").Value = 1400
$ws.Range("L99").Value = 334738
$ws.Range("M99").Value = 98
$ws.Range("N99").Value = -337734
$ws.Range("H105").Value = 1020.9
$ws.Range("I105").Value = 659.8570999999999
$ws.Range("J105").Value = 1863.3334
$ws.Range("K105").Value = 659.8570999999999
$ws.Range("L105").Value = 1863.3334
$ws.Range("M105").Value = 1087.1429
$ws.Range("N105").Value = -5357.3334
$ws.Range("H126").Value = 144259.14
$ws.Range("I126").Value = 1400
$ws.Range("J126").Value = 334738
$ws.Range("K126").Value = 4200
$ws.Range("L126").Value = 1004214
$ws.Range("M126").Value = -1730
$ws.Range("N126").Value = -1009154

# ----- Sheet: CUL -----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 16129670
$ws.Range("I5").Value = 628
$ws.Range("J5").Value = 71429240
$ws.Range("K5").Value = 1884
$ws.Range("L5").Value = 214287720
$ws.Range("M5").Value = -1772
$ws.Range("N5").Value = -214287944
$ws.Range("H34").Value = 41668852
$ws.Range("I34").Value = 160.4
$ws.Range("J34").Value = 71432200
$ws.Range("K34").Value = 481.2
$ws.Range("L34").Value = 214296600
$ws.Range("M34").Value = -397.2
$ws.Range("N34").Value = -214296768
$ws.Range("H135").Value = 16129670
$ws.Range("I135").Value = 628
$ws.Range("J135").Value = 71429240
$ws.Range("K135").Value = 5652
$ws.Range("L135").Value = 642863160
$ws.Range("M135").Value = -3117
$ws.Range("N135").Value = -642868230

# ----- Sheet: GSM -----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 1930.8552
$ws.Range("J5").Value = 1943.2667
$ws.Range("L5").Value = 1943.2667
$ws.Range("N5").Value = -2167.2667
$ws.Range("H70").Value = 4917
$ws.Range("I70").Value = 4832.9165
$ws.Range("K70").Value = 4832.9165
$ws.Range("M70").Value = -4562.9165
$ws.Range("H73").Value = 4917
$ws.Range("I73").Value = 4832.9165
$ws.Range("K73").Value = 4832.9165
$ws.Range("M73").Value = -3896.9165

# ----- Sheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1858.4286
$ws.Range("I7").Value = 2012
$ws.Range("J7").Value = 1582
$ws.Range("K7").Value = 2012
$ws.Range("L7").Value = 1582
$ws.Range("M7").Value = -1900
$ws.Range("N7").Value = -1806
$ws.Range("H22").Value = 355.1111
$ws.Range("I22").Value = 300
$ws.Range("J22").Value = 424
$ws.Range("K22").Value = 300
$ws.Range("L22").Value = 424
$ws.Range("M22").Value = -5
$ws.Range("N22").Value = -1014
$ws.Range("H27").Value = 355.1111
$ws.Range("I27").Value = 300
$ws.Range("J27").Value = 424
$ws.Range("K27").Value = 300
$ws.Range("L27").Value = 424
$ws.Range("M27").Value = -193
$ws.Range("N27").Value = -638
$ws.Range("H119").Value = 0
$ws.Range("J119").Value = 0
$ws.Range("L119").Value = 0
$ws.Range("N119").ClearContents()
$ws.Range("H122").Value = 1928.4324
$ws.Range("I122").Value = 1781.6786
$ws.Range("J122").Value = 2385
$ws.Range("K122").Value = 5345.0358
$ws.Range("L122").Value = 7155
$ws.Range("M122").Value = -2895.0358
$ws.Range("N122").Value = -12055
$ws.Range("H126").Value = 1858.4286
$ws.Range("I126").Value = 2012
$ws.Range("J126").Value = 1582
$ws.Range("K126").Value = 6036
$ws.Range("L126").Value = 4746
$ws.Range("M126").Value = -3566
$ws.Range("N126").Value = -9686
$ws.Range("H128").Value = 40130
$ws.Range("J128").Value = 40130
$ws.Range("L128").Value = 40130
$ws.Range("N128").Value = -50090
$ws.Range("H132").Value = 11853.214
$ws.Range("I132").Value = 14594.8
$ws.Range("J132").Value = 4999.25
$ws.Range("K132").Value = 43784.39999999999
$ws.Range("L132").Value = 14997.75
$ws.Range("M132").Value = -41254.39999999999
$ws.Range("N132").Value = -20057.75
$ws.Range("H133").Value = 26663
$ws.Range("J133").Value = 26663
$ws.Range("L133").Value = 26663
$ws.Range("N133").Value = -31723

# ----- Sheet: WVR -----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2755.932
$ws.Range("I132").Value = 3347.7036
$ws.Range("J132").Value = 1816.0588
$ws.Range("K132").Value = 10043.1108
$ws.Range("L132").Value = 5448.1764
$ws.Range("M132").Value = -7513.110799999999
$ws.Range("N132").Value = -10508.1764
